$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9778419338421145
$ws.Range("C2").Value = 0.3436140131745731
$ws.Range("E2").Value = 0.1337365669603231
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.2984173100465242
$ws.Range("H2").Value = 0.4864538644292793
$ws.Range("I2").Value = 0.4732350977511377
$ws.Range("L2").Value = 0.1978996258906704
$ws.Range("M2").Value = 0.2054740588579449
$ws.Range("O2").Value = 1.497060148339727
$ws.Range("B3").Value = 0.8649347154653242
$ws.Range("C3").Value = 0.3302306120373828
$ws.Range("E3").Value = 0.1356021692935691
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.3041019625657597
$ws.Range("H3").Value = 0.4931147252279686
$ws.Range("I3").Value = 0.4832625443421161
$ws.Range("L3").Value = 0.1953551288554181
$ws.Range("M3").Value = 0.1877879147840247
$ws.Range("O3").Value = 1.52281720800616
$ws.Range("B4").Value = 0.7953819212791018
$ws.Range("C4").Value = 0.322024241992267
$ws.Range("E4").Value = 0.1368155021723507
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.3079307752659481
$ws.Range("H4").Value = 0.4974931133764571
$ws.Range("I4").Value = 0.4897932788413542
$ws.Range("L4").Value = 0.1938916434008462
$ws.Range("M4").Value = 0.1769263345977663
$ws.Range("O4").Value = 1.539945681353494
$ws.Range("B5").Value = 0.7669836037212008
$ws.Range("C5").Value = 0.3186832204373502
$ws.Range("E5").Value = 0.1373270195442335
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3095759265809761
$ws.Range("H5").Value = 0.4993498888062504
$ws.Range("I5").Value = 0.4925484998674534
$ws.Range("L5").Value = 0.1933201879857194
$ws.Range("M5").Value = 0.1724999300994909
$ws.Range("O5").Value = 1.547255459033892
$ws.Range("B6").Value = 0.7622648275964252
$ws.Range("C6").Value = 0.3181286464433413
$ws.Range("E6").Value = 0.1374129882773061
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3098542227841818
$ws.Range("H6").Value = 0.4996625859389958
$ws.Range("I6").Value = 0.4930116698816107
$ws.Range("L6").Value = 0.1932268058623947
$ws.Range("M6").Value = 0.1717649254645579
$ws.Range("O6").Value = 1.548489145190821
$ws.Range("B7").Value = 0.7949991514020098
$ws.Range("C7").Value = 0.3219791706209492
$ws.Range("E7").Value = 0.1368223315221971
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3079526189687627
$ws.Range("H7").Value = 0.4975178607901753
$ws.Range("I7").Value = 0.4898300566917939
$ws.Range("L7").Value = 0.19388383553067
$ws.Range("M7").Value = 0.1768666390138236
$ws.Range("O7").Value = 1.540042928963828
$ws.Range("B8").Value = 0.9389599649086904
$ws.Range("C8").Value = 0.3389973200951601
$ws.Range("E8").Value = 0.1343657539409806
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.300307010441351
$ws.Range("H8").Value = 0.4886906190211704
$ws.Range("I8").Value = 0.4766149045411243
$ws.Range("L8").Value = 0.197001809194802
$ws.Range("M8").Value = 0.19937653121886
$ws.Range("O8").Value = 1.505668252559033
$ws.Range("B9").Value = 1.219381552169693
$ws.Range("C9").Value = 0.3724441692888547
$ws.Range("E9").Value = 0.1300860503091297
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.2880088363446163
$ws.Range("H9").Value = 0.4736711353935519
$ws.Range("I9").Value = 0.453670906615641
$ws.Range("L9").Value = 0.2038978600116153
$ws.Range("M9").Value = 0.2434882737164372
$ws.Range("O9").Value = 1.448704487685774
$ws.Range("B10").Value = 1.424168613769325
$ws.Range("C10").Value = 0.3970470147060041
$ws.Range("E10").Value = 0.127268487784642
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.2806294261029549
$ws.Range("H10").Value = 0.4640337761706803
$ws.Range("I10").Value = 0.4386308268088914
$ws.Range("L10").Value = 0.2094382670979797
$ws.Range("M10").Value = 0.2758652770320893
$ws.Range("O10").Value = 1.413251068095917
$ws.Range("B11").Value = 1.517045039992468
$ws.Range("C11").Value = 0.4082427128059862
$ws.Range("E11").Value = 0.1260574328506388
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.2776347468691753
$ws.Range("H11").Value = 0.459953140546645
$ws.Range("I11").Value = 0.4321844271723778
$ws.Range("L11").Value = 0.212061123179808
$ws.Range("M11").Value = 0.290584838433908
$ws.Range("O11").Value = 1.398518076075703
$ws.Range("B12").Value = 1.552172448023896
$ws.Range("C12").Value = 0.4124824551820154
$ws.Range("E12").Value = 0.1256089852429465
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2765530618178573
$ws.Range("H12").Value = 0.4584515711137627
$ws.Range("I12").Value = 0.4298003192429549
$ws.Range("L12").Value = 0.2130690093367633
$ws.Range("M12").Value = 0.2961571822725944
$ws.Range("O12").Value = 1.393140215254675
$ws.Range("B13").Value = 1.544609075599851
$ws.Range("C13").Value = 0.4115693487958083
$ws.Range("E13").Value = 0.1257051150684689
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.276783690575563
$ws.Range("H13").Value = 0.4587730176480918
$ws.Range("I13").Value = 0.4303112419918067
$ws.Range("L13").Value = 0.2128512918177563
$ws.Range("M13").Value = 0.2949571566526856
$ws.Range("O13").Value = 1.394289475530243
$ws.Range("B14").Value = 1.519935865432956
$ws.Range("C14").Value = 0.4085915175366495
$ws.Range("E14").Value = 0.1260203353764108
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.2775447057288645
$ws.Range("H14").Value = 0.4598287298917469
$ws.Range("I14").Value = 0.431987141912904
$ws.Range("L14").Value = 0.2121437490316538
$ws.Range("M14").Value = 0.2910433133004489
$ws.Range("O14").Value = 1.398071600322496
$ws.Range("B15").Value = 1.504817159873994
$ws.Range("C15").Value = 0.4067675218964837
$ws.Range("E15").Value = 0.1262147388079665
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.278017672243962
$ws.Range("H15").Value = 0.4604810736644112
$ws.Range("I15").Value = 0.4330211068512275
$ws.Range("L15").Value = 0.2117122666625875
$ws.Range("M15").Value = 0.2886457478668376
$ws.Range("O15").Value = 1.400414481683953
$ws.Range("B16").Value = 1.418093045151863
$ws.Range("C16").Value = 0.3963153936393553
$ws.Range("E16").Value = 0.1273490542954163
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.2808324457311784
$ws.Range("H16").Value = 0.4643065651550629
$ws.Range("I16").Value = 0.4390600831710314
$ws.Range("L16").Value = 0.2092689136791961
$ws.Range("M16").Value = 0.2749031096121186
$ws.Range("O16").Value = 1.414242026073524
$ws.Range("B17").Value = 1.36481668542541
$ws.Range("C17").Value = 0.389904049987905
$ws.Range("E17").Value = 0.128063013261474
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.2826521733492839
$ws.Range("H17").Value = 0.4667311313704232
$ws.Range("I17").Value = 0.4428661760080415
$ws.Range("L17").Value = 0.2077961968846154
$ws.Range("M17").Value = 0.2664699135995505
$ws.Range("O17").Value = 1.423082492746929
$ws.Range("B18").Value = 1.334147103199143
$ws.Range("C18").Value = 0.3862167882792562
$ws.Range("E18").Value = 0.1284803157789265
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.2837329149227301
$ws.Range("H18").Value = 0.4681542370632528
$ws.Range("I18").Value = 0.4450925490137125
$ws.Range("L18").Value = 0.2069587809251345
$ws.Range("M18").Value = 0.261618545690979
$ws.Range("O18").Value = 1.428298558017289
$ws.Range("B19").Value = 1.323758440543202
$ws.Range("C19").Value = 0.3849684201048547
$ws.Range("E19").Value = 0.1286227500511666
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.2841046816770074
$ws.Range("H19").Value = 0.4686409799668283
$ws.Range("I19").Value = 0.4458527469430598
$ws.Range("L19").Value = 0.2066769061513867
$ws.Range("M19").Value = 0.259975826466686
$ws.Range("O19").Value = 1.430087153864392
$ws.Range("B20").Value = 1.370490796857723
$ws.Range("C20").Value = 0.390586512440791
$ws.Range("E20").Value = 0.1279863227010454
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.2824549310691111
$ws.Range("H20").Value = 0.4664700760585561
$ws.Range("I20").Value = 0.4424571588802628
$ws.Range("L20").Value = 0.2079519716956497
$ws.Range("M20").Value = 0.2673677288046576
$ws.Range("O20").Value = 1.422127820337622
$ws.Range("B21").Value = 1.527184167269183
$ws.Range("C21").Value = 0.4094661769075287
$ws.Range("E21").Value = 0.1259274720667
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.2773197548422672
$ws.Range("H21").Value = 0.4595174558432831
$ws.Range("I21").Value = 0.4314933412681929
$ws.Range("L21").Value = 0.2123511739974191
$ws.Range("M21").Value = 0.292192951146518
$ws.Range("O21").Value = 1.396955233131493
$ws.Range("B22").Value = 1.629341358559998
$ws.Range("C22").Value = 0.4218059795319675
$ws.Range("E22").Value = 0.1246410670186077
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2742687804027781
$ws.Range("H22").Value = 0.4552281295218137
$ws.Range("I22").Value = 0.4246601806163293
$ws.Range("L22").Value = 0.2153117671258116
$ws.Range("M22").Value = 0.3084079780632507
$ws.Range("O22").Value = 1.381676540988963
$ws.Range("B23").Value = 1.574841861892821
$ws.Range("C23").Value = 0.4152200361043015
$ws.Range("E23").Value = 0.1253222343808686
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.2758691406874618
$ws.Range("H23").Value = 0.4574941123555902
$ws.Range("I23").Value = 0.4282767114038171
$ws.Range("L23").Value = 0.2137238468744016
$ws.Range("M23").Value = 0.2997547164712699
$ws.Range("O23").Value = 1.389723531484591
$ws.Range("B24").Value = 1.367925656413831
$ws.Range("C24").Value = 0.390277975170676
$ws.Range("E24").Value = 0.1280209732031483
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.2825439966840548
$ws.Range("H24").Value = 0.4665880082429936
$ws.Range("I24").Value = 0.4426419565462618
$ws.Range("L24").Value = 0.2078815170282837
$ws.Range("M24").Value = 0.2669618359032384
$ws.Range("O24").Value = 1.422559011741996
$ws.Range("B25").Value = 1.143731208782356
$ws.Range("C25").Value = 0.3633896057731079
$ws.Range("E25").Value = 0.1311863716562504
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.2910459451825886
$ws.Range("H25").Value = 0.4774889656576136
$ws.Range("I25").Value = 0.4595592177628696
$ws.Range("L25").Value = 0.2019488843164794
$ws.Range("M25").Value = 0.2315595087594389
$ws.Range("O25").Value = 1.462993344366865
